$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 through 6 (years 2005-2009), shifting 2010/2011/2012 rows up
$ws.Range("A2:H6").EntireRow.Delete()
